# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# The account-statement table (rows 16-21) is re-sorted in descending
# order by "Periodo Mora" (column E), and the "Salario Basico" amounts
# (column F) follow their row so the 67200 value now lands on period
# 2311 instead of 2306.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2311", "2310", "2309", "2308", "2307", "2306")
$salaries = @(67200, 72000, 72000, 72000, 72000, 72000)

for ($i = 0; $i -lt 6; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $salaries[$i]
}
